$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column FE, shifting FE:FP (and the rest) three
# columns to the right (FE->FH ... FP->FS). This makes room for a new
# "properties.Ngày 18.select.id/name/color" column group, mirroring the other
# per-day select expansions already present on the sheet.
$ws.Range("FE1:FG1").EntireColumn.Insert(-4161)

# New header labels for the inserted columns (row 1).
$ws.Range("FE1").Value = "properties.Ngày 18.select.id"
$ws.Range("FF1").Value = "properties.Ngày 18.select.name"
$ws.Range("FG1").Value = "properties.Ngày 18.select.color"

# Populate the new "Ngày 18" select data for the rows that have attendance
# recorded (everything else stays blank, matching the sheet's existing
# pattern for days without data).
$rows = 8,12,13,17,18
foreach ($r in $rows) {
    $ws.Cells.Item($r, 161).Value = "DjwF"
    $ws.Cells.Item($r, 162).Value = "Đầy đủ"
    $ws.Cells.Item($r, 163).Value = "pink"
}

# Refresh last_edited_time (column D) for the rows whose "Ngày 18" data was
# just added.
foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = "2024-07-18T10:57:00.000Z"
}
